# Code changes to Home suite
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Home_Page): Done / N
$ws.Range("B4").Value = "Done"
$ws.Range("C4").Value = "N"

# Row 11 (ClaimMojio): In Progress / Y
$ws.Range("B11").Value = "In Progress"
$ws.Range("C11").Value = "Y"

# Row 13 (Contact_Support): clear Runmode/Description value
$ws.Range("B13").Value = ""

# Row 17 (TripHistory): clear the "Date search remains" note
$ws.Range("B17").Value = ""

# Update the active selection to A19, matching the saved view state
$ws.Range("A19").Select()
